# Insert a new "fees" column into the currency_movements sheet, between the
# existing "amount" and "currency" columns, so wire transfers can optionally
# carry fees just like the other transaction sheets.
#
# Before: date | buy_date | amount | currency | comment
# After:  date | buy_date | amount | fees     | currency | comment

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currency_movements")

# Shift the "currency" and "comment" columns one place to the right,
# carrying their values/styles with them, and open up column D.
$ws.Range("D:D").Insert()

# Header for the newly inserted column.
$ws.Range("D1").Value = "fees"

# Existing rows default to a fee of 0, written as plain (unstyled) numbers.
$ws.Range("D2:D5").Value = 0
$ws.Range("D2:D5").ClearFormats()
